$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Value = "ffflvdfv"
$ws.Range("F9").Value = "svasvf"
$ws.Range("D7").Value = "asvasv"
$ws.Range("I6").Value = "svasfrfasf"
$ws.Range("G16").Value = "dsfsdfsferq"

$ws.Range("I15").Select()
